$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (refreshed snapshot).
# D-column (Price) values that parse as pure numbers get a temporary
# text NumberFormat so Excel does not silently coerce them to numeric
# cells (the source data is plain text, e.g. "1.00", "0.0000289").

$ws.Range('D2').Value = '96.445.76'
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').Value = '3.699.22'
$ws.Range('E3').Value = '  +2.52%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.74%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.89'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.49%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '654.74'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.71%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.431'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.56%  '

$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.07'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.12%  '

$ws.Range('D11').Value = '3.697.50'
$ws.Range('E11').Value = '  +2.46%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '44.86'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('E13').Value = '  +0.70%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.82'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +5.42%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000289'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +12.21%  '

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '4.389.53'
$ws.Range('E16').Value = '  +2.73%  '

$ws.Range('D17').Value = '96.338.55'
$ws.Range('E17').Value = '  -0.11%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.701.48'
$ws.Range('E18').Value = '  +2.66%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.86'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.51%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.74'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.33%  '

$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.80'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -11.32%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.509'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.66%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '520.57'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.47%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.43'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000209'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.05%  '

$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.02'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.68%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '101.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.31%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '13.26'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.09%  '

$ws.Range('E29').Value = '  -1.97%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.04'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.31%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '12.24'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.55%  '

$ws.Range('E32').Value = '  +0.12%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.88'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.86%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.185'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.33%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '672.86'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +9.12%  '

$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.22%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '32.40'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.61%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.593'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.44%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.85'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.02%  '

$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.84'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +9.52%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.90'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +19.04%  '

$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.161'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.88%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.75%  '

$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.961'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.82%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0451'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.40%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.434'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.35%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.29'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.42%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '23.59'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.55'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.18%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.54'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.30%  '
